$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 273
$ws.Range("A273").Value = 272
$ws.Range("B273").Value = 'Time for some truth'
$ws.Range("C273").Value = 'EmeraldSoul'
$ws.Range("D273").Value = 'Corrupted Bishop'
$ws.Range("F273").Value = 'Round and round we go, The man and his white coat, The key and it''s maker, The chained man, Meals and stories, The ritual for the gates'
$ws.Range("G273").Value = 15
$ws.Range("L273").Value = 100000
$ws.Range("N273").Value = 1000000
$ws.Range("O273").Value = 1000000
$ws.Range("P273").Value = 1000000000
$ws.Range("Q273").Value = 1000
$ws.Range("T273").Value = 1
$ws.Range("AA273").Value = 'It is time to get some answers from this construct. Enough is enough. The visions, the running around, the circles, the backstory—it''s all leading to the same thing: you must be the first one to open the Gates of Heaven and banish the Corrupted Bishop back to where he came from.<br /> <br /> You head to the graveyard where you have been conversing with the construct and search with fevered determination. Nowhere. Nothing. You call out, you scream, you yell—the rage boils in your blood.<br /> <br /> You see the light appear. You hear the small voice—that of a child.<br /> <br /> "What''s wrong? I am here."<br /> <br /> You lay into the construct. Everything. All the details. But then you hear the construct start to cry.<br /> <br /> "This is not true. The Prince created me. I am not a child. I am not something created by anyone but the Prince. These are lies! All lies!"<br /> <br /> You explain that these are not lies, and that the construct needs to face the truth—to help you with the ritual. You demand to know how to open the Gates of Heaven.<br /> <br /> "I don’t know. Honestly. None of this makes sense."<br /> <br /> You tell him who his mother is—the sister of the Madam. How she wept for him when she saw the blue blanket. You demand that something, anything you’ve said, must trigger a memory. All the things you’ve done, all the items you’ve fetched—something must awaken a distant memory.<br /> <br /> "Stop! Stop all of this. None of it''s true. I wasn’t created by priests. I was created by him. The Emerald Prince. I was... I was..."<br /> <br /> He vanishes.<br /> <br /> What do you do now?'
$ws.Range("AB273").Value = 'You walk around the graveyard, fuming with anger. You don’t know what to do. You don’t even know who to talk to. Everything you’ve been told leads you to believe the Construct is the key to the ritual.<br /> <br /> The Prophet and the Church and the Scholar of Hell are both dead. Both lost their lives—one to his own selfish desires to open the gates, likely in connection with the Federation Chancellor and his own ambitions. The other died trying to help you. Alas, his sacrifice was not in vain. You’ve learned all you could from the River Man. There’s nothing left but the Construct itself.<br /> <br /> It appears again, sobbing and crying.<br /> <br /> You approach the construct.<br /> <br /> "You are mean," it says.<br /> <br /> You shake it off.<br /> <br /> "But... you are right," it continues.<br /> <br /> You look surprised.<br /> <br /> "I was hiding my creation even from myself. It wasn’t the Prince—it was the corrupted magics of the Bishop and his priests. I was stolen away, experimented on, tortured, neglected, beaten, and more. All so they could open the gates, gaze upon the face of God, and unleash the angels upon the land. Unleash them all.<br /> <br /> "But only one managed to get through. One who became nameless. One who manipulated the Poet, the Wandering Merchant, and even the Prince himself."<br /> <br /> It pauses. You suddenly know who the construct speaks of.<br /> <br /> You ask if it means the River Man.<br /> <br /> "Yes," the construct replies. After a moment, it continues.<br /> <br /> "He left Heaven. He manipulated the Poet, the Merchant, and the Prince, all in an effort to help the Child. But it was all a long game. He worked with the Sister—my mother—to ensure the Child failed in his ascension, to keep him here.<br /> <br /> "Alas, this caused the Shattering. It ruptured the planes into reality. That wasn’t part of his plan. So he had to help clean up his mess. But when the Child escaped purgatory—that''s when he knew. As he became the Creator, as you know him now, he could see. I’m not sure how, but he could see.<br /> <br /> "He imprisoned the River Man. But you... you freed him. And now he is somewhere out there..."<br /> <br /> You stand, shocked at these revelations. But you focus on what must be done and ask how to start the ritual to banish the Corrupted Bishop.<br /> <br /> "You need the chains, for one. You need the blood of innocence. And you need a relic of the Church. But most of all, you need the incantation—the one that opens the gates."'
$ws.Range("AC273").Value = 'Corrupted Bishop'
$ws.Range("AE273").Value = 7

# Row 274
$ws.Range("A274").Value = 273
$ws.Range("B274").Value = 'The Ritual Ingredients'
$ws.Range("C274").Value = 'EmeraldSoul'
$ws.Range("D274").Value = 'Corrupted Bishop'
$ws.Range("F274").Value = 'Chains of the past, The shadows in The Federation city, The book and the Dagger, A witches incantations'
$ws.Range("H274").Value = 'Heavens Incantation'
$ws.Range("N274").Value = 100000
$ws.Range("O274").Value = 100000
$ws.Range("P274").Value = 1000000000
$ws.Range("Q274").Value = 1000
$ws.Range("U274").Value = 'Time for some truth'
$ws.Range("AA274").Value = 'You ask the construct how to obtain the necessary items. You explain that from what you understand, the blood of innocence refers to the blood of a child—but you don’t feel right about killing a child.<br /> <br /> "I can give that to you—my blood."<br /> <br /> You ask if it will work, pointing out that the construct is not a human child.<br /> <br /> "Yes, it will still work because it came from me as a child. Alas, I need you to secure the rest of the materials first before I can give you the blood. I would start where the River Man was chained up—see if you can retrieve the chains."<br /> <br /> He vanishes. You pause for a moment, then set out. You now know what you need: the items required for a ritual to open the gates and banish the Bishop once and for all.'
$ws.Range("AB274").Value = 'You return with everything you need. The Construct stands before you and asks if you have the incantation. You hand it over and confirm that you have gathered the rest of the items.<br /> <br /> "We can start the ritual now, child. We will need to draw him out—draw him here, to this cemetery."<br /> <br /> You ask about the blood he promised.<br /> <br /> "I have it, child. But it is the last thing we add. We must do this in a very particular order. Are you ready?"<br /> <br /> You nod. Judgment Day has come for the Bishop.'
$ws.Range("AC274").Value = 'Corrupted Bishop'
$ws.Range("AE274").Value = 7

# Row 275
$ws.Range("A275").Value = 274
$ws.Range("B275").Value = 'Chains of the past'
$ws.Range("C275").Value = 'Soldierofthegates'
$ws.Range("D275").Value = 'Corrupted Bishop'
$ws.Range("G275").Value = 20
$ws.Range("M275").Value = 'Chains of The River Man'
$ws.Range("N275").Value = 100000
$ws.Range("O275").Value = 100000
$ws.Range("P275").Value = 1000000000
$ws.Range("Q275").Value = 1250
$ws.Range("U275").Value = 'Time for some truth'
$ws.Range("AA275").Value = 'You arrive back in Purgatory, at the Cellar—but alas, it’s empty. Even the chains are gone. There’s no one here.<br /> <br /> You leave the Cellar and look around, only to see a man standing not far from the entrance. You cautiously approach him. He greets you.<br /> <br /> "Hello, child."<br /> <br /> He looks slightly older than you—fit, dressed in a tunic and leather pants. He carries a sword and a lantern.<br /> <br /> You’ve never seen him before and ask who he is.<br /> <br /> "I am a Soldier of the Gates. I seek the man who was in this Cellar. Do you know anything about him?"<br /> <br /> You begin to explain the ritual, the Key Maker, and how you asked the River Man about it—but you’re cut off by an angry outburst.<br /> <br /> "You did what?! Are you stupid? Do you not know what you’ve done?!"<br /> <br /> You try to explain what the Construct told you and how you planned to fix the mess after dealing with the Bishop, but the Soldier of the Gates isn’t having it.<br /> <br /> "Enough, child! You have done irrevocable harm here today with your idiocy! You stuck your nose where it doesn’t belong and unleashed a madman—an escapee from Heaven itself. Do you know what this means?!"<br /> <br /> You sarcastically state that you’re a thousand steps away from discovering the truth of the Child—but also one step closer.<br /> <br /> "YOU BAFFLING IDIOT!" he yells.<br /> <br /> He starts to storm off. You call after him, asking if he knows where the River Man’s chains are because you need them. He stops and turns.<br /> <br /> "What the hell would you need those for?"<br /> <br /> You begin to tell him how you plan to stop the Corrupted Bishop.'
$ws.Range("AB275").Value = 'After a moment, his rage subsides—if only slightly. A glint of hatred still burns behind his piercing green eyes.<br /> <br /> "You sit here and tell me a story of visions, the Creator, constructs, the Emerald Prince, the Poet, the Wandering Merchant, and the River Man—how it''s all connected to the Shattering. How the sister of the Madam was involved. How all of this ties together. And you expect me to help you by handing over these chains? Meanwhile, you also admit to setting free an Angel of the Lord—who played a role in all of this and was the sole cause of the Shattering?<br /> <br /> "And you expect me to help you?"<br /> <br /> You tell him that in exchange for the chains, you will help him hunt down the River Man. What he does with him is his own business. But you need the chains to at least set one part of this mess right. You warn him that failure to help you could unleash, through the hands of either the Bishop or the Federation Chancellor, the wrath of Heaven—as the Bishop ascends to become the new God.<br /> <br /> He laughs.<br /> <br /> "Oh, you think that’s how it happens? A challenge is declared, and he just wins? He would need the magic of a pure innocent..." His words trail off.<br /> <br /> A pure innocent? you ask.<br /> <br /> "Child... the Construct! He would need the Construct. It’s pure and innocent—created by love, crafted in love."<br /> <br /> You ask how that helps you.<br /> <br /> "You need to stop him from sacrificing the Construct."<br /> <br /> You tell him you can do that—if you get the required items. One of them is the chains.<br /> <br /> "You want me to hand over the chains? Do you realize what will be required of you? The sacrifice?"<br /> <br /> You assure him you have it all under control.<br /> <br /> "And how will you help me in this process?"<br /> <br /> You explain that once you deal with the Bishop, you will return and help him hunt down the River Man.<br /> <br /> "I suppose a deal is struck, then. But mark my words, child. Should you fail—I will be the one you deal with when you approach the gates. And I will not be merciful. Should you break your end of the bargain, I will be the <br /> one to hunt you down myself—and the hell I will unleash..."<br /> <br /> He hands you the chains from his bag.<br /> <br /> You ask one last question: where can you find a relic of the Church, and where can you find the incantations?<br /> <br /> "The relic would come from the Federation. As for the incantation, I would suggest asking the Scholar—but he is dead. So your last bet is the Witch."'
$ws.Range("AC275").Value = 'Corrupted Bishop'
$ws.Range("AE275").Value = 7
$ws.Range("AF275").Value = 'Soldierofthegates'
$ws.Range("AG275").Value = 3

# Row 276
$ws.Range("A276").Value = 275
$ws.Range("B276").Value = 'The shadows in The Federation city'
$ws.Range("C276").Value = 'Shadow'
$ws.Range("D276").Value = 'Corrupted Bishop'
$ws.Range("H276").Value = 'Thieves Tool Kit'
$ws.Range("M276").Value = 'Holy Bejeweled Censer'
$ws.Range("N276").Value = 100000
$ws.Range("O276").Value = 100000
$ws.Range("P276").Value = 1000000000
$ws.Range("Q276").Value = 1350
$ws.Range("U276").Value = 'Chains of the past'
$ws.Range("V276").Value = 'Mothers Sleep'
$ws.Range("AA276").Value = 'You stand within the Federation City. It has been weeks of travel. You wonder to yourself if the Bishop has made any moves—if he has started his preparations. You’ve heard of strange things happening at the cemetery, but you haven’t had time to investigate, as you’ve been focused on searching for the items needed to complete your ritual.<br /> <br /> he visions have been dormant—not haunting you, not plaguing you—but you know that as you get closer and closer to finalizing the details, he, the Bishop himself, will make himself known to you.<br /> <br /> You walk through the streets, keeping to the back alleys, afraid of being apprehended by the Federation soldiers who patrol the city. You study them, along with the people passing by. Now and then, a priest or a Holy Knight walks past. Everyone is dead—these are just distant memories of another time.<br /> <br /> You see a shadow move. You stop and ready yourself for battle, but a soft voice comes from the darkness.<br /> <br /> “Hush now, child. We don’t want to draw more attention to ourselves.”<br /> <br /> You peer into the shadows but see nothing.<br /> <br /> “I am known as Shadow. I hide in the darkness and steal what I please. I can also be hired to assassinate those who stand in the way of others’ political and religious ambitions—for a price, of course. A steep price.”<br /> <br /> You ask if he can be hired to steal something—something you need: a religious artifact.<br /> <br /> “I can. But what would you want with such a toy?”<br /> <br /> You tell him you''re not at liberty to discuss the details, but you can offer a trade in return.<br /> <br /> “What could you possibly offer me that I cannot obtain myself, child?”<br /> <br /> You explain how you can take on things beyond his comprehension—things that can see through the shadows, things that haunt your dreams and torment your waking hours. You describe your unique set of skills for <br /> dealing with what lurks in the depths.<br /> <br /> “Like your visions, child? Visions of what is... what was... visions of deluded memories?”<br /> <br /> You try to mask your shock, but he sees right through it.<br /> <br /> “I know of you, child. I know of your exploits across this and other planes of existence. I know the trouble you’ve caused, the evils you’ve unleashed, and the foes you’ve vanquished. I will help you, because you are like <br /> me—seeking something beyond the shadows.<br /> <br /> “You find me the Thieves’ Toolkit and the Mother’s Sleep, and I will get you the relic you desire.”<br /> <br /> You agree and head off through the alleyways, pondering how this man knows so much about you, yet never reveals himself. Has he always been there, hiding in the shadows?'
$ws.Range("AB276").Value = 'You return to the same alley, searching the shadows for the man known as Shadow.<br /> <br /> You call out in a whisper, careful not to be heard by those passing in the nearby street or peering from the windows above.<br /> <br /> “Child,” comes a voice from behind you, deep in the shadows.<br /> <br /> “I have what you need. Do you have what I want?”<br /> <br /> You hand the items over to the darkness. Slender, cloaked arms reach out and take them from you. In exchange, he places a trinket in your hand—an old censer.<br /> <br /> “I know what you seek to do, child. I have followed you through the depths. I know whom you have spoken with, the pact and the deal you have made, and what you seek to banish. But your quest will be in vain, unless you heed my words: there are delusions, and there is reality. What you seek to destroy will only cause that reality to emerge even more clearly. To make the Creator angry is to face the wrath of the cosmics.”<br /> <br /> With that, he is gone. The shadows are empty. The streets are silent. You look around—but little do you know, you are being watched.'
$ws.Range("AC276").Value = 'Corrupted Bishop'
$ws.Range("AE276").Value = 7
$ws.Range("AF276").Value = 'Shadow'
$ws.Range("AG276").Value = 3

# Row 277
$ws.Range("A277").Value = 276
$ws.Range("B277").Value = 'The book and the Dagger'
$ws.Range("C277").Value = 'MadamTrixi'
$ws.Range("D277").Value = 'Corrupted Bishop'
$ws.Range("H277").Value = 'Old Leather Bound Diary of a sister'
$ws.Range("M277").Value = 'Bone White Dagger'
$ws.Range("N277").Value = 150000
$ws.Range("O277").Value = 150000
$ws.Range("P277").Value = 1500000000
$ws.Range("Q277").Value = 1500
$ws.Range("U277").Value = 'The shadows in The Federation city'
$ws.Range("AA277").Value = 'You start to head out, back toward where you were going. The last thing you need is the incantations—the Witch might have what you''re looking for.<br /> <br /> Alas, as you leave the alleyway, you see priests walking toward you. You turn and pull your cloak hood up, only to see Holy Knights coming from the opposite direction. You turn again and find the Federation Chancellor standing before you.<br /> <br /> "Hello, child! I see you are out and about, consorting with villains and godless heathens, collecting supplies for what I can only assume is a ritual of some kind. Perhaps you''re planning to open the Gates of Heaven? Would you like to peer upon the eyes of God?"<br /> <br /> The knights and priests begin to circle around you, weapons drawn. You prepare yourself for battle and start to conjure your magic—but feel it being drained from you, fading. You feel powerless.<br /> <br /> Looking up, you see the Chancellor holding a device—one that is siphoning your power away.<br /> <br /> "Child, we can''t have you starting a ruckus in our streets now, can we? Kill ''em, boys!"<br /> <br /> You ready your weapons, but as you do, you blink—and suddenly, you are standing before her: the Madam.<br /> <br /> "What fresh hell have you gotten yourself into now, child?"<br /> <br /> You stare at her and ask how you got here. What happened to the Chancellor?<br /> <br /> "Oh, I can put you back—if that’s what you’d like."<br /> <br /> You tell her you’re grateful for her help but curious how she’s able to reach across planes, especially since this one is—or was—locked away.<br /> <br /> "Once you opened the gates to my world, you unlocked my ability to reach across the planes. While I cannot physically leave this one due to the magics that bind me, trust me when I say, child—I will find a way. Now, why don’t you tell me why those nasty old ghosts of the past are after you?"<br /> <br /> You explain that you''re gathering the items required to open the Gates of Heaven and banish the Bishop—and, hopefully, the deluded reality he keeps trying to convince you is real.<br /> <br /> "Oh, but isn’t it?" she asks, laughing.<br /> <br /> You shake it off and tell her you were planning to head to the Witch to ask for the final item—the incantation.<br /> <br /> "Oh, but to get to her, and to get what you want from her, you’ll need something very special. Something I can give you—if you fetch me something very special first."<br /> <br /> You ask what she wants.<br /> <br /> "A simple book. An old one, of sorts. A diary. You fetch this for me, and I’ll give you something for the Witch. A win-win all around."'
$ws.Range("AB277").Value = 'You never thought to ask what the book was for, or what was written in it. Now, you stand aboard a wrecked ship, staring at the book and flipping through its pages—maybe you can glean something from it.<br /> <br /> That’s when you hear her voice:<br /> <br /> "Oh child, you shouldn’t be reading other people’s personal entries. It’s not nice to snoop."<br /> <br /> And before you know it, the pages go blank—and you’re standing before her again.<br /> <br /> She snatches the book from your hands and tosses it into the fire. You ask why, but she only shrugs and walks over to a rotting dresser, barely standing upright, with gaping holes throughout. She opens the second drawer from the top and pulls out a dagger—bone white.<br /> <br /> "A deal is a deal," she says, handing you the dagger.<br /> <br /> "Off to see the Witch. The wonderful Witch of the labyrinth!" she cackles, and then vanishes.<br /> <br /> What was in the book? Why did she burn it?<br /> <br /> Questions for another time. Maybe it wasn’t important to your task... or maybe it was.<br /> <br /> Alas, now you must find your way to the Witch—for the second-to-last item.'
$ws.Range("AC277").Value = 'Corrupted Bishop'
$ws.Range("AE277").Value = 7

# Row 278
$ws.Range("A278").Value = 277
$ws.Range("B278").Value = 'A witches incantations'
$ws.Range("C278").Value = 'TheWitch'
$ws.Range("D278").Value = 'Corrupted Bishop'
$ws.Range("H278").Value = 'Bone White Dagger'
$ws.Range("M278").Value = 'Heavens Incantation'
$ws.Range("N278").Value = 150000
$ws.Range("O278").Value = 150000
$ws.Range("P278").Value = 1500000000
$ws.Range("Q278").Value = 1500
$ws.Range("U278").Value = 'The book and the Dagger'
$ws.Range("V278").Value = 'Vile of Innocence'
$ws.Range("AA278").Value = 'You start to head back toward your destination. The last thing you need is the incantation—the Witch might have what you’re looking for.<br /> <br /> But as you leave the alleyway, you spot priests walking toward you. You turn and pull your cloak’s hood up, only to see Holy Knights coming from the other direction. You spin again—and find the Federation Chancellor standing right in front of you.<br /> <br /> "Hello, child! I see you’re out and about, consorting with villains and godless heathens, collecting supplies for what I can only assume is a ritual of some kind. Perhaps you’re planning to open the Gates of Heaven? Would you like to peer upon the eyes of God?"<br /> <br /> The knights and priests begin circling you, weapons drawn. You ready yourself for battle and begin to summon your magic—but it’s no use. You feel your power being drained, fading. You look up and see the Chancellor holding a device—one that’s pulling the magic from you.<br /> <br /> "Child, we can’t have you starting a ruckus in our streets now, can we? Kill ''em, boys!"<br /> <br /> You draw your weapons, but before you can react, you blink—and suddenly you’re standing before her: the Madam.<br /> <br /> "What fresh hell have you gotten yourself into now, child?"<br /> <br /> You stare at her in confusion and ask how you got here. What about the Chancellor?<br /> <br /> "Oh, I can put you back, if that’s what you’d like."<br /> <br /> You tell her you’re grateful for her help, but you’re curious—how is she able to reach across planes, especially since this one was, or is, locked away?<br /> <br /> "Once you opened the gates to my world, you unlocked my ability to reach across the planes. While I cannot physically leave this plane due to the magics that bind me—trust me when I say, child—I will find a way.<br /> <br /> "Now, why don’t you tell me why those nasty old ghosts of the past are after you?"<br /> <br /> You explain that you’re gathering the items needed to open the Gates of Heaven and banish the Bishop—and, hopefully, the deluded reality he keeps trying to convince you is real.<br /> <br /> "Oh, but isn’t it?" she says with a laugh.<br /> <br /> You shake it off and tell her you were planning to visit the Witch for the final item—the incantation.<br /> <br /> "Oh, but to get to her—and to get what you want from her—you’ll need something very special. Something I can give you... if you fetch me something just as special in return."<br /> <br /> You ask what she wants.<br /> <br /> "A simple book. An old one, a diary. You fetch this for me, and I’ll give you something the Witch will want. A win-win all around."'
$ws.Range("AB278").Value = 'You never thought to ask what the book was for—or what was written inside. Now you stand aboard a wrecked ship, staring down at it. You flip through the pages, hoping to glean something of value.<br /> <br /> That’s when you hear her voice.<br /> <br /> "Oh, child. You shouldn’t be reading other people’s personal entries. It’s not nice to snoop."<br /> <br /> Before you can respond, the pages go blank—and you find yourself standing before her once again.<br /> <br /> She snatches the book from your hands and tosses it into the fire. You ask her why, but she only shrugs and turns away. She walks over to a rotting dresser, barely holding itself together, full of holes and decay. She opens the second drawer from the top and pulls out a dagger—a bone-white dagger.<br /> <br /> "A deal is a deal," she says, handing you the weapon.<br /> <br /> "Off to see the Witch—the wonderful Witch of the labyrinth," she cackles, then vanishes.<br /> <br /> What was in the book? Why did she burn it?<br /> <br /> Questions you’ll have to figure out later. Maybe it wasn’t important to the task at hand... or maybe it was.<br /> <br /> Either way, your path is clear: it’s time to find the Witch, and secure the second-to-last item.'
$ws.Range("AC278").Value = 'Corrupted Bishop'
$ws.Range("AE278").Value = 7
